# "Task ready for kid piloting"
# Replicate the existing stimulus row (row 2) down through row 11 so the
# sheet holds 10 trial rows instead of 1, then leave the last new row
# selected (matching Excel's post-fill-down selection behaviour).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$source = $ws.Range("A2:E2")

for ($row = 3; $row -le 11; $row++) {
    $target = $ws.Range("A" + $row + ":E" + $row)
    $source.Copy($target) | Out-Null
}

$ws.Range("A11:E11").Select() | Out-Null
